$d = $word.ActiveDocument

# Locate the paragraph that currently holds the greeting/prayer sentence
# (identified by distinctive text from the "before" version) rather than
# assuming it is paragraph 1, so the script is resilient to minor
# structural differences.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Today is Sunday*church*") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    $target = $d.Paragraphs.Item(1)
}

$pRange = $target.Range
$contentStart = $pRange.Start
$contentEnd = $pRange.End - 1   # exclude the trailing paragraph mark

$body = $d.Range($contentStart, $contentEnd)

# Rebuild the paragraph's runs from scratch via InsertXML so the exact
# run-split the diff calls for (5 separate <w:r> elements) is preserved
# instead of being re-merged by identical-formatting normalization.
$newXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r><w:t>We are going to church</w:t></w:r>
            <w:r><w:t xml:space="preserve">. God is good and He will answer all my prayers in </w:t></w:r>
            <w:r><w:t>Jesus’</w:t></w:r>
            <w:r><w:t xml:space="preserve"> name.</w:t></w:r>
            <w:r><w:t xml:space="preserve"> I love and fear God, He is very powerful.</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$body.InsertXML($newXml)

Write-Output $d.Content.Text
